$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header for the main image URL column from "mainImageUrl" to "mainImageURI"
$ws.Range("C1").Value = "mainImageURI"

# Move the active cell/selection as recorded in the saved file
$ws.Range("E12").Select()
